# VerifyCatalogSearch.xlsx edit
# Commit: "Moved Test Files to new Package "SearchFunctionalty""
#
# Content changes performed:
#  - Sheet "VerifyCatalogSearch", cell C2: "COMPUTER DESKTOP" -> "DESKTOPs"
#  - Sheet "VerifyCatalogSearch" selection moved from J7 to C6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VerifyCatalogSearch")

# Update the test data value in C2 (was "COMPUTER DESKTOP", now "DESKTOPs")
$ws.Range("C2").Value = "DESKTOPs"

# Move / restore the active selection to C6 as recorded in the saved view state
$ws.Activate()
$ws.Range("C6").Select()
